# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header / footer timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 16:52"

# --- Swap country labels (ranking changed positions) ---
$ws.Range("A68").Value = "Uzbekistan"
$ws.Range("A69").Value = "Oman"

$ws.Range("A136").Value = "Birmania"
$ws.Range("A137").Value = "Gibraltar"

# --- Numeric data updates ---
# Row 4: Estados Unidos
$ws.Range("B4").Value = 850116
$ws.Range("C4").Value = 1399
$ws.Range("E4").Value = 718321
$ws.Range("G4").Value = 78
$ws.Range("H4").Value = 47737

# Row 8: Alemania
$ws.Range("B8").Value = 151175
$ws.Range("C8").Value = 527
$ws.Range("E8").Value = 42521

# Row 18: Suiza
$ws.Range("E18").Value = 7058
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = 1538

# Row 26: Arabia Saudita
$ws.Range("F26").Value = 93

# Row 61: Grecia
$ws.Range("E61").Value = 1706
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 125

# Row 68: Uzbekistan (was Oman's row)
$ws.Range("B68").Value = 1735
$ws.Range("C68").Value = 19
$ws.Range("D68").Value = 503
$ws.Range("E68").Value = 1225
$ws.Range("F68").Value = 8
$ws.Range("H68").Value = 7

# Row 69: Oman (was Uzbekistan's row)
$ws.Range("C69").Value = 102
$ws.Range("D69").Value = 307
$ws.Range("E69").Value = 1401
$ws.Range("F69").Value = 3
$ws.Range("H69").Value = 8

# Row 84: Bulgaria
$ws.Range("B84").Value = 1097
$ws.Range("C84").Value = 73
$ws.Range("E84").Value = 855
$ws.Range("G84").Value = 3
$ws.Range("H84").Value = 52

# Row 95: Libano
$ws.Range("F95").Value = 46

# Row 114: Sri Lanka
$ws.Range("B114").Value = 337
$ws.Range("C114").Value = 7
$ws.Range("E114").Value = 223

# Row 136: Birmania (was Gibraltar's row)
$ws.Range("C136").Value = 9
$ws.Range("D136").Value = 9
$ws.Range("E136").Value = 118
$ws.Range("H136").Value = 5

# Row 137: Gibraltar (was Birmania's row)
$ws.Range("B137").Value = 132
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 127
$ws.Range("E137").Value = 5
$ws.Range("H137").Value = 0

# Row 141: Trinidad y Tobago
$ws.Range("D141").Value = 41
$ws.Range("E141").Value = 66
